# Update the "Date Placeholder 1" shape on every slide that has one:
#  - set the text from 8/6/2019 to 8/9/2019
#  - make sure its position/size is explicit (matches the inherited
#    layout position: 628650,6356351 / 2057400,365125 EMU == 49.5,500.5001 / 162,28.75 pt)
#  - move it to the end of the shape/z-order stack (as PowerPoint does
#    when a placeholder is removed and then re-added via Insert > Header & Footer)

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -eq "Date Placeholder 1") {
            $sh.Left = 49.5
            $sh.Top = 500.5001
            $sh.Width = 162
            $sh.Height = 28.75
            $sh.TextFrame.TextRange.Text = "8/9/2019"
            $sh.ZOrder(0)
        }
    }
}
